$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-51 contain the crypto list (columns B=Coin, C=Link, D=Price, E=Volume(1h)).
# This update refreshes price/volume figures, and re-orders a few rows whose
# underlying data changed position (WrappedEther/Polygon swap, and the
# Aptos/RenderToken/RocketPoolETH/SynthetixNetwork block at the bottom which
# now ends with a new entry "EnergySwap").
#
# Column D values are stored as plain text in the workbook (e.g. "29.957.14",
# "104.00", "0.000007850") even though many of them look numeric. A leading
# apostrophe forces Excel to keep them as literal text instead of silently
# coercing to a number (which would drop trailing zeros / reformat decimals).

$rows = @(
    @{ Row = 2;  B = "Bitcoin";                      C = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc";                                D = "29.957.14";     E = "  +0.26%  " },
    @{ Row = 3;  B = "Ethereum";                      C = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth";                                D = "1.876.56";      E = "  -0.65%  " },
    @{ Row = 4;  B = "TetherUSD";                     C = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt";                              D = "1.001";         E = "  +0.05%  " },
    @{ Row = 5;  B = "XRP";                           C = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp";                                     D = "0.7402";        E = "  -3.94%  " },
    @{ Row = 6;  B = "BNB";                           C = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb";                                     D = "242.71";        E = "  -0.01%  " },
    @{ Row = 7;  B = "USDC";                          C = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc";                                   D = "1.002";         E = "  +0.13%  " },
    @{ Row = 8;  B = "Cardano";                       C = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada";                                 D = "0.3146";        E = "  +0.81%  " },
    @{ Row = 9;  B = "Dogecoin";                      C = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge";                                D = "0.07215";       E = "  +0.48%  " },
    @{ Row = 10; B = "Solana";                        C = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol";                                      D = "24.59";         E = "  -3.96%  " },
    @{ Row = 11; B = "TRON";                          C = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx";                                    D = "0.08321";       E = "  -3.42%  " },
    @{ Row = 12; B = "Polygon";                       C = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic";                                D = "0.7525";        E = "  -1.52%  " },
    @{ Row = 13; B = "WrappedEther";                  C = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth";                           D = "1.893.26";      E = "  -2.44%  " },
    @{ Row = 14; B = "Polkadot";                      C = "https://coinranking.com/coin/25W7FG7om+polkadot-dot";                                    D = "5.417";         E = "  +0.88%  " },
    @{ Row = 15; B = "Litecoin";                      C = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc";                                D = "92.64";         E = "  -1.07%  " },
    @{ Row = 16; B = "WrappedBTC";                    C = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc";                               D = "29.967.13";     E = "  -0.03%  " },
    @{ Row = 17; B = "Uniswap";                       C = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni";                                     D = "6.108";         E = "  -1.24%  " },
    @{ Row = 18; B = "BitcoinCash";                   C = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch";                             D = "250.79";        E = "  +2.54%  " },
    @{ Row = 19; B = "Avalanche";                     C = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax";                                  D = "13.57";         E = "  -1.50%  " },
    @{ Row = 20; B = "ShibaInu";                      C = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib";                                   D = "0.000007850";   E = "  +0.42%  " },
    @{ Row = 21; B = "Dai";                           C = "https://coinranking.com/coin/MoTuySvg7+dai-dai";                                         D = "1.001";         E = "  +0.06%  " },
    @{ Row = 22; B = "WrappedliquidstakedEther2.0";   C = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth";                D = "2.141.51";      E = "  -3.74%  " },
    @{ Row = 23; B = "Chainlink";                     C = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link";                               D = "8.023";         E = "  -0.40%  " },
    @{ Row = 24; B = "BinanceUSD";                    C = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd";                              D = "1.001";         E = "  -0.09%  " },
    @{ Row = 25; B = "Stellar";                       C = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm";                                  D = "0.1551";        E = "  -6.58%  " },
    @{ Row = 26; B = "Cosmos";                        C = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom";                                  D = "9.265";         E = "  -1.12%  " },
    @{ Row = 27; B = "Monero";                        C = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr";                                   D = "165.18";        E = "  +1.75%  " },
    @{ Row = 28; B = "EthereumClassic";                C = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc";                          D = "18.69";         E = "  -0.41%  " },
    @{ Row = 29; B = "LidoDAOToken";                  C = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo";                                 D = "2.035";         E = "  -0.28%  " },
    @{ Row = 30; B = "Toncoin";                       C = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton";                                      D = "1.508";         E = "  +3.70%  " },
    @{ Row = 31; B = "Filecoin";                      C = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil";                                     D = "4.614";         E = "  +2.35%  " },
    @{ Row = 32; B = "PancakeSwap";                   C = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake";                                 D = "1.536";         E = "  +0.13%  " },
    @{ Row = 33; B = "InternetComputer(DFINITY)";     C = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp";                       D = "4.311";         E = "  +5.06%  " },
    @{ Row = 34; B = "Hedera";                        C = "https://coinranking.com/coin/jad286TjB+hedera-hbar";                                      D = "0.05333";       E = "  -2.10%  " },
    @{ Row = 35; B = "ARBITRUM";                      C = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb";                                     D = "1.235";         E = "  -0.39%  " },
    @{ Row = 36; B = "ImmutableX";                    C = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx";                                   D = "0.7477";        E = "  +0.35%  " },
    @{ Row = 37; B = "Frax";                          C = "https://coinranking.com/coin/KfWtaeV1W+frax-frax";                                        D = "1.004";         E = "  +0.08%  " },
    @{ Row = 38; B = "HuobiToken";                    C = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht";                                D = "2.706";         E = "  +0.29%  " },
    @{ Row = 39; B = "VeChain";                       C = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet";                                  D = "0.01974";       E = "  +0.81%  " },
    @{ Row = 40; B = "MXToken";                       C = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx";                                    D = "2.757";         E = "  -0.96%  " },
    @{ Row = 41; B = "TheSandbox";                    C = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand";                                  D = "0.4556";        E = "  +1.99%  " },
    @{ Row = 42; B = "Maker";                         C = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr";                                    D = "1.113.08";      E = "  +0.47%  " },
    @{ Row = 43; B = "FraxShare";                     C = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs";                                    D = "6.133";         E = "  +0.85%  " },
    @{ Row = 44; B = "Aave";                          C = "https://coinranking.com/coin/ixgUfzmLR+aave-aave";                                        D = "72.32";         E = "  -1.15%  " },
    @{ Row = 45; B = "TrustWalletToken";               C = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt";                             D = "0.8558";        E = "  +0.49%  " },
    @{ Row = 46; B = "PaxDollar";                     C = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp";                                   D = "1.003";         E = "  +0.25%  " },
    @{ Row = 47; B = "Quant";                         C = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt";                                    D = "104.00";        E = "  +1.50%  " },
    @{ Row = 48; B = "RenderToken";                   C = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr";                             D = "1.854";         E = "  -0.73%  " },
    @{ Row = 49; B = "Aptos";                         C = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt";                                        D = "7.609";         E = "  -0.68%  " },
    @{ Row = 50; B = "EnergySwap";                    C = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens";                                   D = "9.488";         E = "  -3.08%  " },
    @{ Row = 51; B = "RocketPoolETH";                 C = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth";                                D = "2.038.49";      E = "  -3.54%  " }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = "'" + $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
}
